# Applies the CacheDoc.docx edit described by the commit diff:
#  - wraps the five bold "function name" terms with proofErr spellStart/spellEnd
#    (and, where the trailing colon was glued onto the name run, splits it
#    into its own run so the proofErr wrap covers only the identifier)
#  - rewrites the "checks to make sure ... max cache size." sentence to add
#    the clause "ensuring the cache does not go over the"
#  - relocates the _GoBack bookmark from the end of the last paragraph into
#    the middle of the "This function deals with checking..." paragraph,
#    splitting the word "address" around it

$d = $word.ActiveDocument
$wNs = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

function Set-ParagraphXml($paragraphIndex, $innerXml) {
    $range = $d.Paragraphs($paragraphIndex).Range
    $xml = "<w:p " + $wNs + ">" + $innerXml + "</w:p>"
    $range.InsertXML($xml)
}

# Paragraph 1: "Iplc_sim_init:" heading -> wrap identifier in proofErr
$p1 = "<w:pPr><w:spacing w:line='240' w:lineRule='auto'/><w:rPr><w:b/></w:rPr></w:pPr>"
$p1 = $p1 + "<w:proofErr w:type='spellStart'/>"
$p1 = $p1 + "<w:r><w:rPr><w:b/></w:rPr><w:t>Iplc_sim_init</w:t></w:r>"
$p1 = $p1 + "<w:proofErr w:type='spellEnd'/>"
$p1 = $p1 + "<w:r><w:rPr><w:b/></w:rPr><w:t>:</w:t></w:r>"
Set-ParagraphXml 1 $p1

# Paragraph 2: body text - split "checks...max cache size." into 3 runs
$p2 = "<w:pPr><w:spacing w:line='240' w:lineRule='auto'/></w:pPr>"
$p2 = $p2 + "<w:r><w:t xml:space='preserve'>This function takes an index, the block size and the associativity of the cache. Before allocating the cache it will perform </w:t></w:r>"
$p2 = $p2 + "<w:r><w:t xml:space='preserve'>checks to make sure the cache meets pre-established specifications like </w:t></w:r>"
$p2 = $p2 + "<w:r><w:t>ensuring the cache does not go over the</w:t></w:r>"
$p2 = $p2 + "<w:r><w:t xml:space='preserve'> max cache size. </w:t></w:r>"
$p2 = $p2 + "<w:r><w:t>It uses this information to allocate memory and dynamically create the cache</w:t></w:r>"
$p2 = $p2 + "<w:r><w:t xml:space='preserve'> with the proper associativity</w:t></w:r>"
$p2 = $p2 + "<w:r><w:t xml:space='preserve'>. </w:t></w:r>"
Set-ParagraphXml 2 $p2

# Paragraph 3: "Iplc_sim_trap_address:" heading -> split colon into its own
# run and wrap the identifier in proofErr
$p3 = "<w:pPr><w:spacing w:line='240' w:lineRule='auto'/><w:rPr><w:b/></w:rPr></w:pPr>"
$p3 = $p3 + "<w:proofErr w:type='spellStart'/>"
$p3 = $p3 + "<w:r><w:rPr><w:b/></w:rPr><w:t>Iplc_sim_trap_address</w:t></w:r>"
$p3 = $p3 + "<w:proofErr w:type='spellEnd'/>"
$p3 = $p3 + "<w:r><w:rPr><w:b/></w:rPr><w:t>:</w:t></w:r>"
Set-ParagraphXml 3 $p3

# Paragraph 4: body text - move the _GoBack bookmark into the middle of the
# paragraph, splitting "address" into "a" + bookmark + "ddress"
$p4 = "<w:pPr><w:spacing w:line='240' w:lineRule='auto'/></w:pPr>"
$p4 = $p4 + "<w:r><w:t>This function deals with checking if a given address is in the cache. It takes into account the given associativity, and looks through the cache data structure. It will update the counter for a hit or a miss. After looking through the appropriate entries for the a</w:t></w:r>"
$p4 = $p4 + "<w:bookmarkStart w:id='0' w:name='_GoBack'/>"
$p4 = $p4 + "<w:bookmarkEnd w:id='0'/>"
$p4 = $p4 + "<w:r><w:t>ddress it will call the appropriate function to deal with a hit or a miss.</w:t></w:r>"
Set-ParagraphXml 4 $p4

# Paragraph 5: "Destroy_cache:" heading -> wrap identifier in proofErr
$p5 = "<w:pPr><w:spacing w:line='240' w:lineRule='auto'/><w:rPr><w:b/></w:rPr></w:pPr>"
$p5 = $p5 + "<w:proofErr w:type='spellStart'/>"
$p5 = $p5 + "<w:r><w:rPr><w:b/></w:rPr><w:t>Destroy_cache</w:t></w:r>"
$p5 = $p5 + "<w:proofErr w:type='spellEnd'/>"
$p5 = $p5 + "<w:r><w:rPr><w:b/></w:rPr><w:t>:</w:t></w:r>"
Set-ParagraphXml 5 $p5

# Paragraph 7: "iplc_sim_LRU_replace_on_miss:" heading -> wrap identifier in proofErr
$p7 = "<w:pPr><w:spacing w:line='240' w:lineRule='auto'/><w:rPr><w:b/></w:rPr></w:pPr>"
$p7 = $p7 + "<w:proofErr w:type='spellStart'/>"
$p7 = $p7 + "<w:r><w:rPr><w:b/></w:rPr><w:t>iplc_sim_LRU_replace_on_miss</w:t></w:r>"
$p7 = $p7 + "<w:proofErr w:type='spellEnd'/>"
$p7 = $p7 + "<w:r><w:rPr><w:b/></w:rPr><w:t>:</w:t></w:r>"
Set-ParagraphXml 7 $p7

# Paragraph 9: "iplc_sim_LRU_update_on_hit:" heading -> wrap identifier in proofErr
$p9 = "<w:pPr><w:spacing w:line='240' w:lineRule='auto'/><w:rPr><w:b/></w:rPr></w:pPr>"
$p9 = $p9 + "<w:proofErr w:type='spellStart'/>"
$p9 = $p9 + "<w:r><w:rPr><w:b/></w:rPr><w:t>iplc_sim_LRU_update_on_hit</w:t></w:r>"
$p9 = $p9 + "<w:proofErr w:type='spellEnd'/>"
$p9 = $p9 + "<w:r><w:rPr><w:b/></w:rPr><w:t>:</w:t></w:r>"
Set-ParagraphXml 9 $p9

# Paragraph 10: body text - remove the trailing _GoBack bookmark (moved to
# paragraph 4 above)
$rsquo = [char]0x2019
$p10 = "<w:pPr><w:spacing w:line='240' w:lineRule='auto'/></w:pPr>"
$p10 = $p10 + "<w:r><w:t xml:space='preserve'>This function takes in a cache index and a given associative entry. It is called if an element has already been determined to be in the cache, and it will update the element" + $rsquo + "s information. It will set it as the MRU, and update the previous MRU information on the cache. </w:t></w:r>"
Set-ParagraphXml 10 $p10
